$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix player name bug: "Kuldeep Yadav" -> "kapil dev" in row 12 (A12)
$ws.Range("A12").Value = "kapil dev"
